$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Semestre ideal:" value EM-9 -> EM-10
$ws.Range("B9").Value = "EM-10"
$ws.Range("C9").Value = "EM-10"

# Rebuild the "Requisitos" list (rows 25..55) with the updated set of prerequisites.
$reqs = @(
  "LOB1003 -  Cálculo I  (Requisito)`n",
  "LOB1004 -  Cálculo II  (Requisito)`n",
  "LOB1008 -  Ciência, Tecnologia e Sociedade  (Requisito)`n",
  "LOB1012 -  Estatística  (Requisito)`n",
  "LOB1018 -  Física I  (Requisito)`n",
  "LOB1019 -  Física II  (Requisito)`n",
  "LOB1036 -  Geometria Analítica  (Requisito)`n",
  "LOB1037 -  Àlgebra Linear  (Requisito)`n",
  "LOB1038 -  Física Experimental I  (Requisito)`n",
  "LOB1039 -  Física Experimental III  (Requisito)`n",
  "LOB1041 -  Física Experimental II  (Requisito)`n",
  "LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito)`n",
  "LOB1046 -  Engenharia do Meio Ambiente  (Requisito)`n",
  "LOB1052 -  Cálculo III  (Requisito)`n",
  "LOB1053 -  Física III  (Requisito)`n",
  "LOM3013 -  Ciência dos Materiais  (Requisito)`n",
  "LOM3018 -  Introdução à Engenharia de Materiais  (Requisito)`n",
  "LOM3037 -  Química Inorgânica  (Requisito)`n",
  "LOM3056 -  Fundamentos de Química Orgânica  (Requisito)`n",
  "LOM3099 -  Estática  (Requisito)`n",
  "LOM3104 -  Projeto Integrado em Engenharia de Materiais I  (Requisito)`n",
  "LOM3105 -  Computação e análise de dados em Engenharia  (Requisito)`n",
  "LOM3204 -  Desenho Técnico e Projeto Assistido por Computador  (Requisito)`n",
  "LOQ4095 -  Química Geral Experimental  (Requisito)`n",
  "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito)`n",
  "LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito)`n",
  "LOQ4246 -  Engenharia da Sustentabilidade  (Requisito)`n"
)

for ($i = 0; $i -lt $reqs.Count; $i++) {
  $row = 25 + $i
  $ws.Cells.Item($row, 2).Value = $reqs[$i]
  $ws.Cells.Item($row, 3).Value = $reqs[$i]
}

# The old list had 31 rows (25..55); the new one only has 27 (25..51) -> drop the tail.
$firstRow = 25 + $reqs.Count
$lastRow = 55
if ($firstRow -le $lastRow) {
  $ws.Range("A" + $firstRow + ":C" + $lastRow).EntireRow.Delete() | Out-Null
}
